# Weekly fruit/vegetable price update for "Cebollín baby".
# Two new weekly observations are inserted into the historical series:
#   - a new row at position 71 (pushing the former rows 71-77 down to 72-78)
#   - a new row at position 79 (pushing the former rows 78-94 down to 80-96)
# All other rows (1-70) are untouched, and the sheet dimension grows from
# A1:R94 to A1:R96 automatically as a result of the two inserted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the first new row at 71 (shifts old 71..94 down to 72..95) ---
$ws.Rows.Item(71).Insert()

# --- Insert the second new row at 79 (shifts rows now at 79..95 down to 80..96) ---
$ws.Rows.Item(79).Insert()

# --- Populate the brand-new row 71 ---
$ws.Cells.Item(71, 1).Value = 1
$ws.Cells.Item(71, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(71, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(71, 4).Value = 44230
$ws.Cells.Item(71, 5).Value = 15
$ws.Cells.Item(71, 6).Value = 100112038
$ws.Cells.Item(71, 7).Value = "Cebollín baby"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Primera"
$ws.Cells.Item(71, 10).Value = 250
$ws.Cells.Item(71, 11).Value = 5500
$ws.Cells.Item(71, 12).Value = 6000
$ws.Cells.Item(71, 13).Value = 5750
$ws.Cells.Item(71, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(71, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(71, 16).Value = 2875
$ws.Cells.Item(71, 17).Value = 2
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# --- Populate the brand-new row 79 ---
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 44782
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = 100112038
$ws.Cells.Item(79, 7).Value = "Cebollín baby"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 250
$ws.Cells.Item(79, 11).Value = 3500
$ws.Cells.Item(79, 12).Value = 4000
$ws.Cells.Item(79, 13).Value = 3750
$ws.Cells.Item(79, 14).Value = "$/paquete 1,5 a 2 kilos"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 1875
$ws.Cells.Item(79, 17).Value = 2
$ws.Cells.Item(79, 18).Value = "Hortaliza"

# --- Make sure the date column keeps the date/time number format used
#     throughout column D (style carried over from Insert, but set
#     explicitly too so it is correct regardless of engine defaults). ---
$ws.Cells.Item(71, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
